$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new year column (O) mirroring the formatting of the existing last
# data column (N), then fill in the 2020 data point.

# O4: new year header (2020), copy style/format from N4
$ws.Range("N4").Copy($ws.Range("O4")) | Out-Null
$ws.Range("O4").Value = 2020

# O5: new data value (83.3), copy style/format from N5
$ws.Range("N5").Copy($ws.Range("O5")) | Out-Null
$ws.Range("O5").Value = 83.3

$excel.CutCopyMode = 0

# Update the selection to match the authored state after the edit.
$ws.Range("O12").Select() | Out-Null

$wb.Save()
